$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated theta_se values (row 4), columns B..L
$ws.Range("B4").Value = "(1.1)"
$ws.Range("C4").Value = "(0.03)"
$ws.Range("D4").Value = "(0.62)"
$ws.Range("E4").Value = "(0.19)"
$ws.Range("F4").Value = "(0.94)"
$ws.Range("G4").Value = "(0.77)"
$ws.Range("H4").Value = "(0.78)"
$ws.Range("I4").Value = "(0.15)"
$ws.Range("J4").Value = "(0.66)"
$ws.Range("K4").Value = "(0.26)"
$ws.Range("L4").Value = "(2.39)"

# Updated lambda_se values (row 6), columns B..L
$ws.Range("B6").Value = "(0.52)"
$ws.Range("C6").Value = "(0.13)"
$ws.Range("D6").Value = "(0.1)"
$ws.Range("E6").Value = "(0.55)"
$ws.Range("F6").Value = "(0.15)"
$ws.Range("G6").Value = "(0.51)"
$ws.Range("H6").Value = "(0.07)"
$ws.Range("I6").Value = "(0.71)"
$ws.Range("J6").Value = "(0.95)"
$ws.Range("K6").Value = "(0.13)"
$ws.Range("L6").Value = "(1.38)"
